# Updated symbol list on Mon Dec 19 11:34:10 UTC 2022 with GitHub Actions
#
# This script reproduces the refreshed crypto price/volume snapshot:
#  - Many "Price" (column D) values are refreshed to newer quotes.
#  - Rows 42/43 (BKEXToken / CEJI) swapped ranking order, so their
#    Coin / Link / Price / Volume columns are exchanged (with a couple
#    of the price+volume figures also refreshed to new values).
#
# All of the touched cells originally hold their numbers as literal text
# (inline strings), so we force the Range to Text format before writing
# the value -- this stops Excel from "helpfully" re-interpreting strings
# like "0.1300" or "0.005040" as numbers (which would silently drop the
# significant trailing zeros or flip into scientific notation for very
# small values). ClearFormats() afterwards removes the temporary "@"
# number-format again so the cell keeps the same (unformatted) look it
# had before, while the stored value remains text.

function Set-TextValue {
    param($ws, $addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D "Price" refreshes (row number -> new text value) ---
$priceChanges = [ordered]@{
    "D2"  = "247.88"
    "D4"  = "5.516"
    "D5"  = "0.05691"
    "D7"  = "0.8081"
    "D8"  = "1.036"
    "D9"  = "0.1502"
    "D10" = "0.08104"
    "D11" = "0.03147"
    "D12" = "0.03022"
    "D13" = "0.09291"
    "D14" = "3.466"
    "D15" = "0.001654"
    "D16" = "0.04708"
    "D17" = "0.0005863"
    "D18" = "0.006350"
    "D19" = "0.005040"
    "D21" = "0.0001502"
    "D22" = "0.0003202"
    "D23" = "3.766"
    "D24" = "6.429"
    "D25" = "2.116"
    "D26" = "0.3315"
    "D27" = "0.1300"
    "D40" = "0.04115"
    "D41" = "0.006969"
    "D44" = "0.008913"
    "D45" = "0.00005883"
    "D47" = "0.0005503"
    "D48" = "0.6828"
    "D49" = "0.008625"
}

foreach ($addr in $priceChanges.Keys) {
    Set-TextValue $ws $addr $priceChanges[$addr]
}

# --- Rows 42 & 43 swap places (BKEXToken <-> CEJI) plus refreshed values ---
Set-TextValue $ws "B42" "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue $ws "D42" "0.003504"
Set-TextValue $ws "E42" "41CEJICEJI"

Set-TextValue $ws "B43" "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue $ws "D43" "0.1044"
Set-TextValue $ws "E43" "42BKEXTokenBKK"
